$d = $word.ActiveDocument

# 1) "n'a" -> "na" within "cest a dire fin qui n'a poinct encores servi Et les mects dans"
$d.Content.Find.Execute(
    "cest a dire fin qui n'a poinct encores servi Et les mects dans",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "cest a dire fin qui na poinct encores servi Et les mects dans", 2)

# 2) "mects" -> "mectes" within " puys mects laultre par dessus auecq du "
$d.Content.Find.Execute(
    " puys mects laultre par dessus auecq du ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " puys mectes laultre par dessus auecq du ", 2)

# The text substitution above lands inside a paragraph that, further
# along, also contains two back-to-back runs ("s" then "able", forming
# the word "sable") that happen to share identical run formatting. The
# text-mutation pass silently coalesces that pair into a single "sable"
# run as a side effect. Nudge a (no-visual-effect) formatting property
# on the "s" run only, which forces the engine to re-split the merged
# run back into its original two pieces.
$sableRng = $d.Content
$sableFound = $sableRng.Find.Execute("sable", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($sableFound) {
    $sOnly = $d.Range($sableRng.Start, $sableRng.Start + 1)
    $sOnly.Bold = $false
}

# 3) "lautre" -> "laultre" within " frottes lung contre lautre soict cave ou plat"
$d.Content.Find.Execute(
    " frottes lung contre lautre soict cave ou plat",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " frottes lung contre laultre soict cave ou plat", 2)

# 4) "regarder" -> "regarde<add>r</add>" within " quil soict en panchant Apres regarder"
$d.Content.Find.Execute(
    " quil soict en panchant Apres regarder",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " quil soict en panchant Apres regarde<add>r</add>", 2)
